# Auto-generated Excel COM-interop script
# Applies numeric cell value updates across multiple worksheets
# as described by the target OOXML diff (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 58 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 48.090908
$ws.Range("I4").Value = 48.090908
$ws.Range("K4").Value = 48.090908
$ws.Range("M4").Value = 65.909092
$ws.Range("H33").Value = 1079952.8
$ws.Range("J33").Value = 5624.75
$ws.Range("L33").Value = 5624.75
$ws.Range("N33").Value = -6082.75
$ws.Range("H58").Value = 985
$ws.Range("J58").Value = 1950
$ws.Range("L58").Value = 5850
$ws.Range("N58").Value = -6150
$ws.Range("H62").Value = 20006.928
$ws.Range("I62").Value = 21000.5
$ws.Range("J62").Value = 19261.75
$ws.Range("K62").Value = 21000.5
$ws.Range("L62").Value = 19261.75
$ws.Range("M62").Value = -20376.5
$ws.Range("N62").Value = -20509.75
$ws.Range("H65").Value = 20006.928
$ws.Range("I65").Value = 21000.5
$ws.Range("J65").Value = 19261.75
$ws.Range("K65").Value = 105002.5
$ws.Range("L65").Value = 96308.75
$ws.Range("M65").Value = -101882.5
$ws.Range("N65").Value = -102548.75
$ws.Range("H70").Value = 45180
$ws.Range("I70").Value = 55225
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 165675
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -165405
$ws.Range("N70").Value = -15540
$ws.Range("H73").Value = 45180
$ws.Range("I73").Value = 55225
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 165675
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -164739
$ws.Range("N73").Value = -16872
$ws.Range("H135").Value = 727.7778
$ws.Range("I135").Value = 625.9091
$ws.Range("J135").Value = 887.8570999999999
$ws.Range("K135").Value = 5633.1819
$ws.Range("L135").Value = 7990.7139
$ws.Range("M135").Value = -3098.1819
$ws.Range("N135").Value = -13060.7139
$ws.Range("H138").Value = 2274.375
$ws.Range("I138").Value = 1245.8667
$ws.Range("J138").Value = 3181.8823
$ws.Range("K138").Value = 3737.6001
$ws.Range("L138").Value = 9545.6469
$ws.Range("M138").Value = 1402.3999
$ws.Range("N138").Value = -19825.6469
$ws.Range("H141").Value = 28699.576
$ws.Range("I141").Value = 28699.576
$ws.Range("K141").Value = 86098.728
$ws.Range("M141").Value = -80918.728

# --- Sheet ARM: 57 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7061564.5
$ws.Range("I32").Value = 1164264.9
$ws.Range("K32").Value = 1164264.9
$ws.Range("M32").Value = -1163977.9
$ws.Range("H61").Value = 2079.6428
$ws.Range("I61").Value = 2159.6667
$ws.Range("K61").Value = 2159.6667
$ws.Range("M61").Value = -1947.6667
$ws.Range("H74").Value = 1462.0667
$ws.Range("I74").Value = 1505.04
$ws.Range("K74").Value = 1505.04
$ws.Range("M74").Value = -631.04
$ws.Range("H75").Value = 27115.334
$ws.Range("J75").Value = 27115.334
$ws.Range("L75").Value = 27115.334
$ws.Range("N75").Value = -28863.334
$ws.Range("H77").Value = 1462.0667
$ws.Range("I77").Value = 1505.04
$ws.Range("K77").Value = 7525.2
$ws.Range("M77").Value = -3157.2
$ws.Range("H78").Value = 27115.334
$ws.Range("J78").Value = 27115.334
$ws.Range("L78").Value = 81346.00199999999
$ws.Range("N78").Value = -90082.00199999999
$ws.Range("H88").Value = 11110.7
$ws.Range("I88").Value = 977
$ws.Range("J88").Value = 17866.5
$ws.Range("K88").Value = 977
$ws.Range("L88").Value = 17866.5
$ws.Range("M88").Value = -571
$ws.Range("N88").Value = -18678.5
$ws.Range("H91").Value = 11110.7
$ws.Range("I91").Value = 977
$ws.Range("J91").Value = 17866.5
$ws.Range("K91").Value = 977
$ws.Range("L91").Value = 17866.5
$ws.Range("M91").Value = 427
$ws.Range("N91").Value = -20674.5
$ws.Range("H122").Value = 6618.694
$ws.Range("I122").Value = 7509.6216
$ws.Range("J122").Value = 3871.6667
$ws.Range("K122").Value = 22528.8648
$ws.Range("L122").Value = 11615.0001
$ws.Range("M122").Value = -20078.8648
$ws.Range("N122").Value = -16515.0001
$ws.Range("H136").Value = 2079.6428
$ws.Range("I136").Value = 2159.6667
$ws.Range("K136").Value = 6479.000100000001
$ws.Range("M136").Value = -3929.000100000001
$ws.Range("H137").Value = 90000
$ws.Range("J137").Value = 90000
$ws.Range("L137").Value = 90000
$ws.Range("N137").Value = -100200
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet BSM: 29 cell update(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 533.3333
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -646
$ws.Range("H94").Value = 1612.7715
$ws.Range("I94").Value = 1450.069
$ws.Range("J94").Value = 2399.1667
$ws.Range("K94").Value = 1450.069
$ws.Range("L94").Value = 2399.1667
$ws.Range("M94").Value = -999.069
$ws.Range("N94").Value = -3301.1667
$ws.Range("H99").Value = 2344.5715
$ws.Range("I99").Value = 2294.1538
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2294.1538
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -796.1538
$ws.Range("N99").Value = -5996
$ws.Range("H107").Value = 3804.3784
$ws.Range("I107").Value = 2364.44
$ws.Range("J107").Value = 6804.25
$ws.Range("K107").Value = 2364.44
$ws.Range("L107").Value = 6804.25
$ws.Range("M107").Value = -444.4400000000001
$ws.Range("N107").Value = -10644.25
$ws.Range("H134").Value = 2864.8628
$ws.Range("I134").Value = 2232.725
$ws.Range("K134").Value = 6698.174999999999
$ws.Range("M134").Value = -4163.174999999999

# --- Sheet CRP: 8 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6035.4287
$ws.Range("I31").Value = 3130.75
$ws.Range("K31").Value = 3130.75
$ws.Range("M31").Value = -2835.75
$ws.Range("H34").Value = 6035.4287
$ws.Range("I34").Value = 3130.75
$ws.Range("K34").Value = 3130.75
$ws.Range("M34").Value = -2928.75

# --- Sheet CUL: 8 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2142.0667
$ws.Range("J68").Value = 2145.7646
$ws.Range("L68").Value = 6437.293799999999
$ws.Range("N68").Value = -8059.293799999999
$ws.Range("H71").Value = 2142.0667
$ws.Range("J71").Value = 2145.7646
$ws.Range("L71").Value = 19311.8814
$ws.Range("N71").Value = -27423.8814

# --- Sheet GSM: 19 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 624.5
$ws.Range("I2").Value = 713.8570999999999
$ws.Range("J2").Value = 499.4
$ws.Range("K2").Value = 713.8570999999999
$ws.Range("L2").Value = 499.4
$ws.Range("M2").Value = -600.8570999999999
$ws.Range("N2").Value = -725.4
$ws.Range("H107").Value = 9653.849
$ws.Range("I107").Value = 13657.4
$ws.Range("K107").Value = 13657.4
$ws.Range("M107").Value = -11737.4
$ws.Range("H122").Value = 3357.5
$ws.Range("I122").Value = 3175
$ws.Range("K122").Value = 9525
$ws.Range("M122").Value = -7075
$ws.Range("H132").Value = 4177.488
$ws.Range("I132").Value = 4163.778
$ws.Range("K132").Value = 12491.334
$ws.Range("M132").Value = -9961.334000000001

# --- Sheet LTW: 15 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H132").Value = 2788.4
$ws.Range("I132").Value = 2593.8
$ws.Range("J132").Value = 3761.4
$ws.Range("K132").Value = 7781.400000000001
$ws.Range("L132").Value = 11284.2
$ws.Range("M132").Value = -5251.400000000001
$ws.Range("N132").Value = -16344.2

# --- Sheet WVR: 26 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 80321.63
$ws.Range("J62").Value = 97578.21000000001
$ws.Range("L62").Value = 97578.21000000001
$ws.Range("N62").Value = -98826.21000000001
$ws.Range("H65").Value = 80321.63
$ws.Range("J65").Value = 97578.21000000001
$ws.Range("L65").Value = 487891.05
$ws.Range("N65").Value = -494131.05
$ws.Range("H96").Value = 1340.6
$ws.Range("I96").Value = 1434.3334
$ws.Range("J96").Value = 1200
$ws.Range("K96").Value = 1434.3334
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = -61.33339999999998
$ws.Range("N96").Value = -3946
$ws.Range("H132").Value = 7239.72
$ws.Range("I132").Value = 10044.059
$ws.Range("J132").Value = 1280.5
$ws.Range("K132").Value = 30132.177
$ws.Range("L132").Value = 3841.5
$ws.Range("M132").Value = -27602.177
$ws.Range("N132").Value = -8901.5
$ws.Range("H136").Value = 3509.8147
$ws.Range("I136").Value = 3610.4
$ws.Range("K136").Value = 10831.2
$ws.Range("M136").Value = -8281.200000000001

Write-Host "Applied all Leviathan Profits updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
